$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.564566135406494
$ws.Range("B1").Value = 4.348228454589844
$ws.Range("C1").Value = 3.164056301116943
$ws.Range("D1").Value = 1.29225480556488
$ws.Range("E1").Value = 0.9297963380813599
